# Lab1Instructions.pptx - slide 1 "Subtitle 2" placeholder edit
#
# Before:
#   Para 1: "Malcolm " + "Machesky" + " and Adrian Kirchner"
#   Para 2: "Aka team \u201cYamaha piano\u201d"
#
# After (single paragraph, 7 runs):
#   "Team Yamaha piano: Malcolm " + "Machesky" + " " + "seatNo"
#   + ": 45 CWID: A20414760, Adrian Kirchner " + "seatNo" + ": 34 CWID: A20425060"

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# ---------------------------------------------------------------------------
# Step 1: merge paragraph 2 into paragraph 1.
#
# Replace everything from just after "Machesky" (the " and Adrian Kirchner"
# tail of paragraph 1) through the end of paragraph 2 with the full new
# tail text, in one assignment. Because the target range crosses the
# paragraph boundary, the engine folds paragraph 2's text into paragraph 1
# - it leaves a now-empty trailing paragraph behind, which we explicitly
# delete afterwards to get back to a single, clean paragraph.
# ---------------------------------------------------------------------------
$tailStart = 17                                  # 1-based offset of " and Adrian Kirchner" in the full range
$tailLen   = $tr.Length - $tailStart + 1         # through the end of paragraph 2 (inclusive of the para break)

$newTail = " seatNo: 45 CWID: A20414760, Adrian Kirchner seatNo: 34 CWID: A20425060"
$tr.Characters($tailStart, $tailLen).Text = $newTail

if ($tr.Paragraphs().Count -gt 1) {
    $tr.Paragraphs(2).Delete()
}

# ---------------------------------------------------------------------------
# Step 2: split the new tail run into the individual runs shown in the diff.
# These edits are length-preserving (same text going back in, just carved
# into smaller pieces) so the offsets computed against $newTail remain valid.
# ---------------------------------------------------------------------------
$pieces = @(
    " ",
    "seatNo",
    ": 45 CWID: A20414760, Adrian Kirchner ",
    "seatNo",
    ": 34 CWID: A20425060"
)

$offset = $tailStart
foreach ($piece in $pieces) {
    $len = $piece.Length
    $tr.Characters($offset, $len).Text = $piece
    $offset += $len
}

# ---------------------------------------------------------------------------
# Step 3: leftmost edit last, so it doesn't disturb the offsets used above.
# "Malcolm " -> "Team Yamaha piano: Malcolm "
# ---------------------------------------------------------------------------
$tr.Characters(1, 8).Text = "Team Yamaha piano: Malcolm "
